$d = $word.ActiveDocument

# Locate the Subtitle paragraph that carries the "Doctors are warning..."
# byline so the edit is anchored to content rather than a hard-coded index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Doctors are warning on dangerous TikTok challenge*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Doctors are warning on dangerous TikTok challenge' paragraph"
}

# Drop the paragraph out of the (style-inherited) numbering/list definition
# by stamping direct numPr/numId=0 formatting on it - same effect as the
# Word UI's "None" list button.
[void]$target.Range.ListFormat.RemoveNumbers()

# "warning on" -> "warning about"
[void]$target.Range.Find.Execute("warning on dangerous", $true, $false, $false, $false, $false, $true, 1, $false, "warning about dangerous", 2)
